$wb = $excel.ActiveWorkbook

# Switch calculation to manual (adds calcMode="manual" to calcPr)
$excel.Calculation = -4135

# --- Inputs sheet: add the source row and the array formula that spills it ---
$wsInputs = $wb.Worksheets.Item("Inputs")
$wsInputs.Range("A5").Value = 1
$wsInputs.Range("B5").Value = 2
$wsInputs.Range("C5").Value = 3
$wsInputs.Range("D5").Value = 4
$wsInputs.Range("E5").Value = 5

$wsInputs.Range("A6:E6").FormulaArray = "=A5:E5"

[void]$wsInputs.Range("A7").Select()

# --- Outputs sheet: label + formula referencing the array result ---
$wsOutputs = $wb.Worksheets.Item("Outputs")
$wsOutputs.Range("A15").Value = "Arraying formula"
$wsOutputs.Range("B15").Formula = "=Inputs!E6"

[void]$wsOutputs.Range("A16").Select()
